$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.458056666666667
$ws.Range("H2").Value = 4.374169999999999
$ws.Range("I2").Value = 0.2323568509805328
$ws.Range("J2").Value = 0.2323568509805327
$ws.Range("M2").Value = 86.89540866666668
$ws.Range("N2").Value = 260.686226
$ws.Range("O2").Value = 0.319779657009892
$ws.Range("P2").Value = 0.3197796570098919
$ws.Range("Q2").Value = 126.6984299091578
$ws.Range("R2").Value = 1140.28586918242
$ws.Range("S2").Value = 0.07430299411045335
$ws.Range("T2").Value = 0.07430299411045332
$ws.Range("G3").Value = 1.458056666666667
$ws.Range("H3").Value = 4.374169999999999
$ws.Range("I3").Value = 0.2323568509805328
$ws.Range("J3").Value = 0.2323568509805327
$ws.Range("O3").Value = 0.1999969065479545
$ws.Range("P3").Value = 0.1999969065479545
$ws.Range("Q3").Value = 79.23985622866111
$ws.Range("R3").Value = 713.1587060579499
$ws.Range("S3").Value = 0.04647065141133061
$ws.Range("T3").Value = 0.0464706514113306
$ws.Range("G4").Value = 1.458056666666667
$ws.Range("H4").Value = 4.374169999999999
$ws.Range("I4").Value = 0.2323568509805328
$ws.Range("J4").Value = 0.2323568509805327
$ws.Range("M4").Value = 60.92601633333334
$ws.Range("N4").Value = 182.778049
$ws.Range("O4").Value = 0.224210932487692
$ws.Range("P4").Value = 0.224210932487692
$ws.Range("Q4").Value = 88.83358428825889
$ws.Range("R4").Value = 799.50225859433
$ws.Range("S4").Value = 0.05209694622824895
$ws.Range("T4").Value = 0.05209694622824893
$ws.Range("G5").Value = 1.458056666666667
$ws.Range("H5").Value = 4.374169999999999
$ws.Range("I5").Value = 0.2323568509805328
$ws.Range("J5").Value = 0.2323568509805327
$ws.Range("M5").Value = 7.809668333333332
$ws.Range("N5").Value = 23.429005
$ws.Range("O5").Value = 0.02873998867505581
$ws.Range("P5").Value = 0.02873998867505581
$ws.Range("Q5").Value = 11.38693897787222
$ws.Range("R5").Value = 102.48245080085
$ws.Range("S5").Value = 0.006677933265752143
$ws.Range("T5").Value = 0.006677933265752142
$ws.Range("G6").Value = 1.458056666666667
$ws.Range("H6").Value = 4.374169999999999
$ws.Range("I6").Value = 0.2323568509805328
$ws.Range("J6").Value = 0.2323568509805327
$ws.Range("M6").Value = 61.75795633333333
$ws.Range("N6").Value = 185.273869
$ws.Range("O6").Value = 0.2272725152794058
$ws.Range("P6").Value = 0.2272725152794058
$ws.Range("Q6").Value = 90.04659995152555
$ws.Range("R6").Value = 810.4193995637298
$ws.Range("S6").Value = 0.05280832596474776
$ws.Range("T6").Value = 0.05280832596474775
$ws.Range("H7").Value = 5.708772
$ws.Range("I7").Value = 0.3032511962008422
$ws.Range("J7").Value = 0.3032511962008422
$ws.Range("M7").Value = 86.89540866666668
$ws.Range("N7").Value = 260.686226
$ws.Range("O7").Value = 0.319779657009892
$ws.Range("P7").Value = 0.3197796570098919
$ws.Range("Q7").Value = 165.355358641608
$ws.Range("R7").Value = 1488.198227774472
$ws.Range("S7").Value = 0.09697356350894479
$ws.Range("T7").Value = 0.09697356350894477
$ws.Range("H8").Value = 5.708772
$ws.Range("I8").Value = 0.3032511962008422
$ws.Range("J8").Value = 0.3032511962008422
$ws.Range("O8").Value = 0.1999969065479545
$ws.Range("P8").Value = 0.1999969065479545
$ws.Range("R8").Value = 930.7503944062199
$ws.Range("S8").Value = 0.06064930114713527
$ws.Range("T8").Value = 0.06064930114713526
$ws.Range("H9").Value = 5.708772
$ws.Range("I9").Value = 0.3032511962008422
$ws.Range("J9").Value = 0.3032511962008422
$ws.Range("M9").Value = 60.92601633333334
$ws.Range("N9").Value = 182.778049
$ws.Range("O9").Value = 0.224210932487692
$ws.Range("P9").Value = 0.224210932487692
$ws.Range("Q9").Value = 115.937578705092
$ws.Range("R9").Value = 1043.438208345828
$ws.Range("S9").Value = 0.06799223347819888
$ws.Range("T9").Value = 0.06799223347819887
$ws.Range("H10").Value = 5.708772
$ws.Range("I10").Value = 0.3032511962008422
$ws.Range("J10").Value = 0.3032511962008422
$ws.Range("M10").Value = 7.809668333333332
$ws.Range("N10").Value = 23.429005
$ws.Range("O10").Value = 0.02873998867505581
$ws.Range("P10").Value = 0.02873998867505581
$ws.Range("Q10").Value = 14.86120530354
$ws.Range("R10").Value = 133.75084773186
$ws.Range("S10").Value = 0.008715435944509334
$ws.Range("T10").Value = 0.008715435944509334
$ws.Range("H11").Value = 5.708772
$ws.Range("I11").Value = 0.3032511962008422
$ws.Range("J11").Value = 0.3032511962008422
$ws.Range("M11").Value = 61.75795633333333
$ws.Range("N11").Value = 185.273869
$ws.Range("O11").Value = 0.2272725152794058
$ws.Range("P11").Value = 0.2272725152794058
$ws.Range("Q11").Value = 117.520697297652
$ws.Range("R11").Value = 1057.686275678868
$ws.Range("S11").Value = 0.06892066212205401
$ws.Range("T11").Value = 0.068920662122054
$ws.Range("G12").Value = 2.914094333333333
$ws.Range("H12").Value = 8.742283
$ws.Range("I12").Value = 0.4643919528186251
$ws.Range("J12").Value = 0.4643919528186251
$ws.Range("M12").Value = 86.89540866666668
$ws.Range("N12").Value = 260.686226
$ws.Range("O12").Value = 0.319779657009892
$ws.Range("P12").Value = 0.3197796570098919
$ws.Range("Q12").Value = 253.2214179882176
$ws.Range("R12").Value = 2278.992761893958
$ws.Range("S12").Value = 0.1485030993904939
$ws.Range("T12").Value = 0.1485030993904938
$ws.Range("G13").Value = 2.914094333333333
$ws.Range("H13").Value = 8.742283
$ws.Range("I13").Value = 0.4643919528186251
$ws.Range("J13").Value = 0.4643919528186251
$ws.Range("O13").Value = 0.1999969065479545
$ws.Range("P13").Value = 0.1999969065479545
$ws.Range("Q13").Value = 158.3699874559672
$ws.Range("R13").Value = 1425.329887103705
$ws.Range("S13").Value = 0.09287695398948867
$ws.Range("T13").Value = 0.09287695398948866
$ws.Range("G14").Value = 2.914094333333333
$ws.Range("H14").Value = 8.742283
$ws.Range("I14").Value = 0.4643919528186251
$ws.Range("J14").Value = 0.4643919528186251
$ws.Range("M14").Value = 60.92601633333334
$ws.Range("N14").Value = 182.778049
$ws.Range("O14").Value = 0.224210932487692
$ws.Range("P14").Value = 0.224210932487692
$ws.Range("Q14").Value = 177.5441589495408
$ws.Range("R14").Value = 1597.897430545867
$ws.Range("S14").Value = 0.1041217527812442
$ws.Range("T14").Value = 0.1041217527812442
$ws.Range("G15").Value = 2.914094333333333
$ws.Range("H15").Value = 8.742283
$ws.Range("I15").Value = 0.4643919528186251
$ws.Range("J15").Value = 0.4643919528186251
$ws.Range("M15").Value = 7.809668333333332
$ws.Range("N15").Value = 23.429005
$ws.Range("O15").Value = 0.02873998867505581
$ws.Range("P15").Value = 0.02873998867505581
$ws.Range("Q15").Value = 22.75811023537944
$ws.Range("R15").Value = 204.822992118415
$ws.Range("S15").Value = 0.01334661946479434
$ws.Range("T15").Value = 0.01334661946479434
$ws.Range("G16").Value = 2.914094333333333
$ws.Range("H16").Value = 8.742283
$ws.Range("I16").Value = 0.4643919528186251
$ws.Range("J16").Value = 0.4643919528186251
$ws.Range("M16").Value = 61.75795633333333
$ws.Range("N16").Value = 185.273869
$ws.Range("O16").Value = 0.2272725152794058
$ws.Range("P16").Value = 0.2272725152794058
$ws.Range("Q16").Value = 179.9685105892141
$ws.Range("R16").Value = 1619.716595302927
$ws.Range("S16").Value = 0.1055435271926041
$ws.Range("T16").Value = 0.1055435271926041
